$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.396551333333333
$ws.Range("H2").Value = 7.189654
$ws.Range("I2").Value = 0.7965668642594712
$ws.Range("J2").Value = 0.7965668642594713
$ws.Range("M2").Value = 0.1808983333333333
$ws.Range("N2").Value = 0.5426949999999999
$ws.Range("O2").Value = 0.09546831801815302
$ws.Range("P2").Value = 0.09546831801815302
$ws.Range("Q2").Value = 0.4335321419477777
$ws.Range("R2").Value = 3.901789277529999
$ws.Range("S2").Value = 0.07604689871984613
$ws.Range("T2").Value = 0.07604689871984613

$ws.Range("G3").Value = 2.396551333333333
$ws.Range("H3").Value = 7.189654
$ws.Range("I3").Value = 0.7965668642594712
$ws.Range("J3").Value = 0.7965668642594713
$ws.Range("M3").Value = 1.572737
$ws.Range("N3").Value = 4.718211
$ws.Range("O3").Value = 0.8300051930177132
$ws.Range("P3").Value = 0.8300051930177132
$ws.Range("Q3").Value = 3.769144954332667
$ws.Range("R3").Value = 33.922304588994
$ws.Range("S3").Value = 0.661154633921197
$ws.Range("T3").Value = 0.6611546339211971

$ws.Range("G4").Value = 2.396551333333333
$ws.Range("H4").Value = 7.189654
$ws.Range("I4").Value = 0.7965668642594712
$ws.Range("J4").Value = 0.7965668642594713
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1412166666666667
$ws.Range("N4").Value = 0.42365
$ws.Range("O4").Value = 0.0745264889641337
$ws.Range("P4").Value = 0.07452648896413369
$ws.Range("Q4").Value = 0.3384329907888889
$ws.Range("R4").Value = 3.0458969171
$ws.Range("S4").Value = 0.05936533161842807
$ws.Range("T4").Value = 0.05936533161842807

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.074599
$ws.Range("H5").Value = 0.223797
$ws.Range("I5").Value = 0.02479525085917582
$ws.Range("J5").Value = 0.02479525085917582
$ws.Range("M5").Value = 0.1808983333333333
$ws.Range("N5").Value = 0.5426949999999999
$ws.Range("O5").Value = 0.09546831801815302
$ws.Range("P5").Value = 0.09546831801815302
$ws.Range("Q5").Value = 0.01349483476833333
$ws.Range("R5").Value = 0.121453512915
$ws.Range("S5").Value = 0.002367160894363679
$ws.Range("T5").Value = 0.002367160894363679

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.074599
$ws.Range("H6").Value = 0.223797
$ws.Range("I6").Value = 0.02479525085917582
$ws.Range("J6").Value = 0.02479525085917582
$ws.Range("M6").Value = 1.572737
$ws.Range("N6").Value = 4.718211
$ws.Range("O6").Value = 0.8300051930177132
$ws.Range("P6").Value = 0.8300051930177132
$ws.Range("Q6").Value = 0.117324607463
$ws.Range("R6").Value = 1.055921467167
$ws.Range("S6").Value = 0.02058018697529285
$ws.Range("T6").Value = 0.02058018697529285

$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.074599
$ws.Range("H7").Value = 0.223797
$ws.Range("I7").Value = 0.02479525085917582
$ws.Range("J7").Value = 0.02479525085917582
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.1412166666666667
$ws.Range("N7").Value = 0.42365
$ws.Range("O7").Value = 0.0745264889641337
$ws.Range("P7").Value = 0.07452648896413369
$ws.Range("Q7").Value = 0.01053462211666667
$ws.Range("R7").Value = 0.09481159905
$ws.Range("S7").Value = 0.001847902989519293
$ws.Range("T7").Value = 0.001847902989519293

$ws.Range("G8").Value = 0.53745
$ws.Range("H8").Value = 1.61235
$ws.Range("I8").Value = 0.1786378848813529
$ws.Range("J8").Value = 0.1786378848813529
$ws.Range("M8").Value = 0.1808983333333333
$ws.Range("N8").Value = 0.5426949999999999
$ws.Range("O8").Value = 0.09546831801815302
$ws.Range("P8").Value = 0.09546831801815302
$ws.Range("Q8").Value = 0.09722380924999997
$ws.Range("R8").Value = 0.8750142832499999
$ws.Range("S8").Value = 0.0170542584039432
$ws.Range("T8").Value = 0.01705425840394321

$ws.Range("G9").Value = 0.53745
$ws.Range("H9").Value = 1.61235
$ws.Range("I9").Value = 0.1786378848813529
$ws.Range("J9").Value = 0.1786378848813529
$ws.Range("M9").Value = 1.572737
$ws.Range("N9").Value = 4.718211
$ws.Range("O9").Value = 0.8300051930177132
$ws.Range("P9").Value = 0.8300051930177132
$ws.Range("Q9").Value = 0.84526750065
$ws.Range("R9").Value = 7.60740750585
$ws.Range("S9").Value = 0.1482703721212233
$ws.Range("T9").Value = 0.1482703721212233

$ws.Range("G10").Value = 0.53745
$ws.Range("H10").Value = 1.61235
$ws.Range("I10").Value = 0.1786378848813529
$ws.Range("J10").Value = 0.1786378848813529
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.1412166666666667
$ws.Range("N10").Value = 0.42365
$ws.Range("O10").Value = 0.0745264889641337
$ws.Range("P10").Value = 0.07452648896413369
$ws.Range("Q10").Value = 0.0758968975
$ws.Range("R10").Value = 0.6830720775
$ws.Range("S10").Value = 0.01331325435618633
$ws.Range("T10").Value = 0.01331325435618633

